# Auto-generated edit script: refresh cached market-price / profit figures
# (columns H:N) on rows scattered across the eight crafting-job sheets.
# Values are plain cached numbers (no formulas), sourced from an external
# market-data refresh per the commit message.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1668.3125
$ws.Range("I9").Value = 1637.6
$ws.Range("J9").Value = 1682.2727
$ws.Range("K9").Value = 1637.6
$ws.Range("L9").Value = 1682.2727
$ws.Range("M9").Value = -1468.6
$ws.Range("N9").Value = -2020.2727
$ws.Range("H19").Value = 792.63635
$ws.Range("I19").Value = 868.3333
$ws.Range("J19").Value = 701.8
$ws.Range("K19").Value = 868.3333
$ws.Range("L19").Value = 701.8
$ws.Range("M19").Value = -693.3333
$ws.Range("N19").Value = -1051.8
$ws.Range("H62").Value = 4360.1665
$ws.Range("I62").Value = 4140.4287
$ws.Range("K62").Value = 4140.4287
$ws.Range("M62").Value = -3516.4287
$ws.Range("H64").Value = 5355.0557
$ws.Range("J64").Value = 5350
$ws.Range("L64").Value = 5350
$ws.Range("N64").Value = -5846
$ws.Range("H65").Value = 4360.1665
$ws.Range("I65").Value = 4140.4287
$ws.Range("K65").Value = 20702.1435
$ws.Range("M65").Value = -17582.1435
$ws.Range("H67").Value = 5355.0557
$ws.Range("J67").Value = 5350
$ws.Range("L67").Value = 5350
$ws.Range("N67").Value = -7066
$ws.Range("H112").Value = 992004.5600000001
$ws.Range("J112").Value = 1211717.2
$ws.Range("L112").Value = 3635151.6
$ws.Range("N112").Value = -3637367.6
$ws.Range("H137").Value = 16803.188
$ws.Range("J137").Value = 53500
$ws.Range("L137").Value = 160500
$ws.Range("N137").Value = -165600
$ws.Range("H138").Value = 6439.1177
$ws.Range("I138").Value = 3037
$ws.Range("J138").Value = 6768.355
$ws.Range("K138").Value = 9111
$ws.Range("L138").Value = 20305.065
$ws.Range("M138").Value = -3971
$ws.Range("N138").Value = -30585.065

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3151.2
$ws.Range("I61").Value = 3001.3333
$ws.Range("K61").Value = 3001.3333
$ws.Range("M61").Value = -2789.3333
$ws.Range("H74").Value = 2166.3872
$ws.Range("I74").Value = 990
$ws.Range("J74").Value = 2726.5715
$ws.Range("K74").Value = 990
$ws.Range("L74").Value = 2726.5715
$ws.Range("M74").Value = -116
$ws.Range("N74").Value = -4474.5715
$ws.Range("H77").Value = 2166.3872
$ws.Range("I77").Value = 990
$ws.Range("J77").Value = 2726.5715
$ws.Range("K77").Value = 4950
$ws.Range("L77").Value = 13632.8575
$ws.Range("M77").Value = -582
$ws.Range("N77").Value = -22368.8575
$ws.Range("H97").Value = 1248
$ws.Range("I97").Value = 1235
$ws.Range("J97").Value = 1495
$ws.Range("K97").Value = 1235
$ws.Range("L97").Value = 1495
$ws.Range("M97").Value = -739
$ws.Range("N97").Value = -2487
$ws.Range("H132").Value = 31252340
$ws.Range("J132").Value = 71430850
$ws.Range("L132").Value = 214292550
$ws.Range("N132").Value = -214297610
$ws.Range("H136").Value = 3151.2
$ws.Range("I136").Value = 3001.3333
$ws.Range("K136").Value = 9003.999899999999
$ws.Range("M136").Value = -6453.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 19423.46
$ws.Range("J100").Value = 19423.46
$ws.Range("L100").Value = 19423.46
$ws.Range("N100").Value = -21587.46
$ws.Range("H134").Value = 38898480
$ws.Range("I134").Value = 20843654
$ws.Range("J134").Value = 111117784
$ws.Range("K134").Value = 62530962
$ws.Range("L134").Value = 333353352
$ws.Range("M134").Value = -62528427
$ws.Range("N134").Value = -333358422

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1281.6666
$ws.Range("I16").Value = 1281.6666
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1281.6666
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -994.6666
$ws.Range("N16").ClearContents()
$ws.Range("H97").Value = 65849.5
$ws.Range("J97").Value = 119999
$ws.Range("L97").Value = 119999
$ws.Range("N97").Value = -121981
$ws.Range("H113").Value = 1281.6666
$ws.Range("I113").Value = 1281.6666
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1281.6666
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 888.3334
$ws.Range("N113").ClearContents()
$ws.Range("H134").Value = 2632982.2
$ws.Range("I134").Value = 1419.871
$ws.Range("K134").Value = 4259.613
$ws.Range("M134").Value = -1724.613
$ws.Range("H141").Value = 303359
$ws.Range("J141").Value = 349865
$ws.Range("L141").Value = 349865
$ws.Range("N141").Value = -360225

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 116.53846
$ws.Range("I33").Value = 64.8
$ws.Range("J33").Value = 148.875
$ws.Range("K33").Value = 388.8
$ws.Range("L33").Value = 893.25
$ws.Range("M33").Value = -105.8
$ws.Range("N33").Value = -1459.25
$ws.Range("H86").Value = 342
$ws.Range("J86").Value = 172
$ws.Range("L86").Value = 516
$ws.Range("N86").Value = -2888
$ws.Range("H89").Value = 342
$ws.Range("J89").Value = 172
$ws.Range("L89").Value = 1548
$ws.Range("N89").Value = -13404
$ws.Range("H110").Value = 2000
$ws.Range("I110").Value = 2000
$ws.Range("K110").Value = 6000
$ws.Range("M110").Value = -1910
$ws.Range("H113").Value = 1407
$ws.Range("J113").Value = 1407
$ws.Range("L113").Value = 4221
$ws.Range("N113").Value = -8561
$ws.Range("H117").Value = 167331
$ws.Range("J117").Value = 1000000
$ws.Range("L117").Value = 3000000
$ws.Range("N117").Value = -3006884
$ws.Range("H131").Value = 663383.5600000001
$ws.Range("I131").Value = 201563.6
$ws.Range("J131").Value = 759596.0600000001
$ws.Range("K131").Value = 604690.8
$ws.Range("L131").Value = 2278788.18
$ws.Range("M131").Value = -599650.8
$ws.Range("N131").Value = -2288868.18
$ws.Range("H132").Value = 5244.1924
$ws.Range("I132").Value = 2595.6667
$ws.Range("J132").Value = 7514.357
$ws.Range("K132").Value = 23361.0003
$ws.Range("L132").Value = 67629.213
$ws.Range("M132").Value = -20831.0003
$ws.Range("N132").Value = -72689.213
$ws.Range("H137").Value = 13785.429
$ws.Range("J137").Value = 15024.5
$ws.Range("L137").Value = 45073.5
$ws.Range("N137").Value = -55273.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 124.052635
$ws.Range("J2").Value = 102
$ws.Range("L2").Value = 102
$ws.Range("N2").Value = -328
$ws.Range("H97").Value = 594.2727
$ws.Range("I97").Value = 571.4286
$ws.Range("J97").Value = 634.25
$ws.Range("K97").Value = 571.4286
$ws.Range("L97").Value = 634.25
$ws.Range("M97").Value = -75.42859999999996
$ws.Range("N97").Value = -1626.25
$ws.Range("H122").Value = 2937
$ws.Range("I122").Value = 3292.9412
$ws.Range("K122").Value = 9878.8236
$ws.Range("M122").Value = -7428.8236

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 21276.615
$ws.Range("I22").Value = 1513.4286
$ws.Range("K22").Value = 1513.4286
$ws.Range("M22").Value = -1218.4286
$ws.Range("H27").Value = 21276.615
$ws.Range("I27").Value = 1513.4286
$ws.Range("K27").Value = 1513.4286
$ws.Range("M27").Value = -1406.4286
$ws.Range("H40").Value = 7011.7144
$ws.Range("I40").Value = 6986.5557
$ws.Range("K40").Value = 6986.5557
$ws.Range("M40").Value = -6850.5557
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").ClearContents()
$ws.Range("H122").Value = 21192
$ws.Range("I122").Value = 13518.45
$ws.Range("J122").Value = 40375.875
$ws.Range("K122").Value = 40555.35000000001
$ws.Range("L122").Value = 121127.625
$ws.Range("M122").Value = -38105.35000000001
$ws.Range("N122").Value = -126027.625
$ws.Range("H136").Value = 55563016
$ws.Range("I136").Value = 7899.353
$ws.Range("J136").Value = 1000000000
$ws.Range("K136").Value = 23698.059
$ws.Range("L136").Value = 3000000000
$ws.Range("M136").Value = -21148.059
$ws.Range("N136").Value = -3000005100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 17159
$ws.Range("I32").Value = 1499
$ws.Range("J32").Value = 24989
$ws.Range("K32").Value = 1499
$ws.Range("L32").Value = 24989
$ws.Range("M32").Value = -1182
$ws.Range("N32").Value = -25623
$ws.Range("H62").Value = 4538.231
$ws.Range("I62").Value = 4727.273
$ws.Range("K62").Value = 4727.273
$ws.Range("M62").Value = -4103.273
$ws.Range("H65").Value = 4538.231
$ws.Range("I65").Value = 4727.273
$ws.Range("K65").Value = 23636.365
$ws.Range("M65").Value = -20516.365
$ws.Range("H81").Value = 11248.75
$ws.Range("I81").Value = 1664.8334
$ws.Range("J81").Value = 16999.1
$ws.Range("K81").Value = 3329.6668
$ws.Range("L81").Value = 33998.2
$ws.Range("M81").Value = -2268.6668
$ws.Range("N81").Value = -36120.2
$ws.Range("H84").Value = 11248.75
$ws.Range("I84").Value = 1664.8334
$ws.Range("J84").Value = 16999.1
$ws.Range("K84").Value = 16648.334
$ws.Range("L84").Value = 169991
$ws.Range("M84").Value = -11344.334
$ws.Range("N84").Value = -180599
$ws.Range("H126").Value = 3541.8333
$ws.Range("I126").Value = 3541.8333
$ws.Range("K126").Value = 10625.4999
$ws.Range("M126").Value = -8155.499899999999
$ws.Range("H136").Value = 1916.1666
$ws.Range("I136").Value = 1799
$ws.Range("J136").Value = 2033.3334
$ws.Range("K136").Value = 5397
$ws.Range("L136").Value = 6100.0002
$ws.Range("M136").Value = -2847
$ws.Range("N136").Value = -11200.0002

